$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '70.412.53'
$ws.Range('E2').Value = '  +0.22%  '
Set-TextValue 'D3' '3.604.96'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '584.33'
$ws.Range('E5').Value = '  -1.27%  '
Set-TextValue 'D6' '190.46'
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('E7').Value = '  -2.01%  '
Set-TextValue 'D8' '3.600.14'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  +0.06%  '
Set-TextValue 'D10' '0.183'
$ws.Range('E10').Value = '  +2.50%  '
Set-TextValue 'D11' '0.666'
$ws.Range('E11').Value = '  +0.30%  '
Set-TextValue 'D12' '56.14'
$ws.Range('E12').Value = '  -4.22%  '
$ws.Range('E13').Value = '  +8.10%  '
Set-TextValue 'D14' '9.75'
$ws.Range('E14').Value = '  -1.49%  '
Set-TextValue 'D15' '4.180.69'
$ws.Range('E15').Value = '  -0.22%  '
Set-TextValue 'D16' '20.02'
$ws.Range('E16').Value = '  +1.59%  '
Set-TextValue 'D17' '3.601.71'
$ws.Range('E17').Value = '  -0.27%  '
Set-TextValue 'D18' '70.339.11'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  +1.03%  '
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  -0.33%  '
Set-TextValue 'D22' '491.47'
$ws.Range('E22').Value = '  +0.25%  '
Set-TextValue 'D23' '20.27'
$ws.Range('E23').Value = '  +4.72%  '
Set-TextValue 'D24' '4.95'
$ws.Range('E24').Value = '  -7.81%  '
Set-TextValue 'D25' '97.11'
$ws.Range('E25').Value = '  +6.71%  '
$ws.Range('E26').Value = '  -1.85%  '
Set-TextValue 'D27' '3.00'
$ws.Range('E27').Value = '  -4.25%  '
Set-TextValue 'D28' '11.16'
$ws.Range('E28').Value = '  -0.62%  '
Set-TextValue 'D29' '9.55'
$ws.Range('E29').Value = '  -0.41%  '
Set-TextValue 'D30' '32.48'
$ws.Range('E30').Value = '  -1.36%  '
Set-TextValue 'D31' '7.67'
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('E32').Value = '  -0.29%  '
Set-TextValue 'D33' '0.120'
$ws.Range('E33').Value = '  +0.70%  '
Set-TextValue 'D34' '66.42'
$ws.Range('E34').Value = '  +0.79%  '
Set-TextValue 'D35' '579.79'
$ws.Range('E35').Value = '  -8.06%  '
Set-TextValue 'D36' '39.11'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('E38').Value = '  +0.07%  '
Set-TextValue 'D39' '0.400'
$ws.Range('E39').Value = '  -1.45%  '
Set-TextValue 'D40' '3.29'
$ws.Range('E40').Value = '  +20.63%  '
Set-TextValue 'D41' '2.92'
$ws.Range('E41').Value = '  +6.62%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D42' '3.47'
$ws.Range('E42').Value = '  -2.35%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D43' '0.137'
$ws.Range('E43').Value = '  -6.22%  '
Set-TextValue 'D44' '3.223.46'
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('E46').Value = '  -1.07%  '
Set-TextValue 'D47' '9.63'
$ws.Range('E47').Value = '  +6.28%  '
Set-TextValue 'D48' '3.36'
$ws.Range('E48').Value = '  +3.07%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('E51').Value = '  -2.85%  '
